# "CC Feature List" status update
#  - New functionality: Toggle walls            -> A20 "Toggle Wall" status goes Inactive -> Complete
#  - Fixed items-near-blocks disappearing bug    -> A5  "Block/Mud"   status goes Inactive -> Complete
#  - Added movement buttons (touchscreen-ready)  -> A16 "Bomb"        status goes Inactive -> Prototype
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "Complete"
$ws.Range("B16").Value = "Prototype"
$ws.Range("B20").Value = "Complete"

# Clear the stray D9 cell selection left over from editing, resetting the view to A1
$ws.Range("A1").Select()
